# making heart and lungs move; making them link to heart and lung pages
# Remove the static shapes that used to compose the "heart" graphic on
# slide 2 (Snip Same Side Corner Rectangle 8, Oval 4, Trapezoid 5).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$namesToRemove = @("Snip Same Side Corner Rectangle 8", "Oval 4", "Trapezoid 5")

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($namesToRemove -contains $shape.Name) {
        $shape.Delete()
    }
}
